$d = $word.ActiveDocument

# The first three paragraphs each had the proper-noun "majankov"/"commitam" split
# into its own run and wrapped in proofErr spell-check markers. Re-typing the full
# sentence via Find/Replace (matching the whole rendered paragraph text) collapses
# those runs back into a single run and drops the now-stale proofErr tags, just
# like a real retype-and-accept-the-squiggle edit in Word would.
$d.Content.Find.Execute("Prvi redak u dokumentu napisala majankov.", $false, $false, $false, $false, $false, $true, 1, $false, "Prvi redak u dokumentu napisala majankov.", 2)
$d.Content.Find.Execute("Drugi redak isto napisala majankov.", $false, $false, $false, $false, $false, $true, 1, $false, "Drugi redak isto napisala majankov.", 2)
$d.Content.Find.Execute("Treći redak isto tako, napravila promjene da ih opet commitam.", $false, $false, $false, $false, $false, $true, 1, $false, "Treći redak isto tako, napravila promjene da ih opet commitam.", 2)

# Append a new paragraph after "COMMIT s malog laptopa ..." with the new commit note.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
# A trailing placeholder "X" is typed first so the bookmark below can be anchored
# at a position strictly before the document's final character (anchoring exactly
# at end-of-document is unreliable here); the placeholder is trimmed off afterwards
# and the bookmark (being anchored to content) stays put.
$newPara.Range.Text = "IDUĆI MERDZ, KONFLIKT 10.3.2020. 5:49 pmX"

$pos = $newPara.Range.End - 2

# Move the "_GoBack" bookmark (left by Word at the site of the last edit) from the
# end of the old last paragraph to the end of this new paragraph.
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos))

$delRange = $d.Range($pos, $pos + 1)
$delRange.Delete()
